$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 797.3
$ws.Range("I12").Value = 254.8
$ws.Range("J12").Value = 1339.8
$ws.Range("K12").Value = 254.8
$ws.Range("L12").Value = 1339.8
$ws.Range("M12").Value = -84.80000000000001
$ws.Range("N12").Value = -1679.8

$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H70").Value = 3470.5625
$ws.Range("I70").Value = 1266.25
$ws.Range("K70").Value = 3798.75
$ws.Range("M70").Value = -3528.75

$ws.Range("H73").Value = 3470.5625
$ws.Range("I73").Value = 1266.25
$ws.Range("K73").Value = 3798.75
$ws.Range("M73").Value = -2862.75

$ws.Range("H113").Value = 6195.88
$ws.Range("I113").Value = 5353.9375
$ws.Range("K113").Value = 5353.9375
$ws.Range("M113").Value = -2099.9375

$ws.Range("H129").Value = 985.6
$ws.Range("I129").Value = 894.8823
$ws.Range("K129").Value = 2684.6469
$ws.Range("M129").Value = 2315.3531

$ws.Range("H135").Value = 2762
$ws.Range("I135").Value = 2643.3333
$ws.Range("J135").Value = 3474
$ws.Range("K135").Value = 23789.9997
$ws.Range("L135").Value = 31266
$ws.Range("M135").Value = -21254.9997
$ws.Range("N135").Value = -36336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 920.5
$ws.Range("I2").Value = 711.0714
$ws.Range("K2").Value = 711.0714
$ws.Range("M2").Value = -598.0714

$ws.Range("H45").Value = 3832
$ws.Range("I45").Value = 3658.4
$ws.Range("J45").Value = 4049
$ws.Range("K45").Value = 3658.4
$ws.Range("L45").Value = 4049
$ws.Range("M45").Value = -3281.4
$ws.Range("N45").Value = -4803

$ws.Range("H74").Value = 2048.8262
$ws.Range("I74").Value = 1807.4
$ws.Range("K74").Value = 1807.4
$ws.Range("M74").Value = -933.4000000000001

$ws.Range("H77").Value = 2048.8262
$ws.Range("I77").Value = 1807.4
$ws.Range("K77").Value = 9037
$ws.Range("M77").Value = -4669

$ws.Range("H116").Value = 920.5
$ws.Range("I116").Value = 711.0714
$ws.Range("K116").Value = 711.0714
$ws.Range("M116").Value = 1582.9286

$ws.Range("H132").Value = 1538.4348
$ws.Range("I132").Value = 1437.909
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 4313.727000000001
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = -1783.727000000001
$ws.Range("N132").Value = -16310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 920.5
$ws.Range("I3").Value = 711.0714
$ws.Range("K3").Value = 711.0714
$ws.Range("M3").Value = -597.0714

$ws.Range("H99").Value = 4033.5833
$ws.Range("I99").Value = 3441.7144
$ws.Range("J99").Value = 4862.2
$ws.Range("K99").Value = 3441.7144
$ws.Range("L99").Value = 4862.2
$ws.Range("M99").Value = -1943.7144
$ws.Range("N99").Value = -7858.2

$ws.Range("H105").Value = 3114.7856
$ws.Range("I105").Value = 3114.7856
$ws.Range("K105").Value = 3114.7856
$ws.Range("M105").Value = -1367.7856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H31").Value = 1996.5714
$ws.Range("I31").Value = 1996.5714
$ws.Range("K31").Value = 1996.5714
$ws.Range("M31").Value = -1701.5714

$ws.Range("H34").Value = 1996.5714
$ws.Range("I34").Value = 1996.5714
$ws.Range("K34").Value = 1996.5714
$ws.Range("M34").Value = -1794.5714

$ws.Range("H99").Value = 1990.4546
$ws.Range("J99").Value = 2027.8572
$ws.Range("L99").Value = 2027.8572
$ws.Range("N99").Value = -5023.8572

$ws.Range("H126").Value = 1990.4546
$ws.Range("J126").Value = 2027.8572
$ws.Range("L126").Value = 6083.571599999999
$ws.Range("N126").Value = -11023.5716

$ws.Range("H132").Value = 2014.1578
$ws.Range("I132").Value = 2014.1578
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6042.4734
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3512.4734
$ws.Range("N132").ClearContents()

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H141").Value = 196177.25
$ws.Range("J141").Value = 212774.14
$ws.Range("L141").Value = 212774.14
$ws.Range("N141").Value = -223134.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2423758.8
$ws.Range("J11").Value = 450
$ws.Range("L11").Value = 1350
$ws.Range("N11").Value = -1630

$ws.Range("H62").Value = 7497.3335
$ws.Range("J62").Value = 10746.5
$ws.Range("L62").Value = 32239.5
$ws.Range("N62").Value = -33611.5

$ws.Range("H65").Value = 7497.3335
$ws.Range("J65").Value = 10746.5
$ws.Range("L65").Value = 96718.5
$ws.Range("N65").Value = -103582.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3895.2222
$ws.Range("I113").Value = 1136.3334
$ws.Range("J113").Value = 5274.6665
$ws.Range("K113").Value = 1136.3334
$ws.Range("L113").Value = 5274.6665
$ws.Range("M113").Value = 1033.6666
$ws.Range("N113").Value = -9614.666499999999

$ws.Range("H132").Value = 1836.9615
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1634.8
$ws.Range("I16").Value = 1671.6923
$ws.Range("J16").Value = 1395
$ws.Range("K16").Value = 1671.6923
$ws.Range("L16").Value = 1395
$ws.Range("M16").Value = -1501.6923
$ws.Range("N16").Value = -1735

$ws.Range("H68").Value = 2305.4119
$ws.Range("I68").Value = 2263.7144
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 2263.7144
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -1514.7144
$ws.Range("N68").Value = -3998

$ws.Range("H71").Value = 2305.4119
$ws.Range("I71").Value = 2263.7144
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 11318.572
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -7574.572
$ws.Range("N71").Value = -19988

$ws.Range("H132").Value = 1995.5438
$ws.Range("I132").Value = 1879.9429
$ws.Range("J132").Value = 2179.4546
$ws.Range("K132").Value = 5639.8287
$ws.Range("L132").Value = 6538.3638
$ws.Range("M132").Value = -3109.8287
$ws.Range("N132").Value = -11598.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1196.9
$ws.Range("I107").Value = 1094.8334
$ws.Range("K107").Value = 3284.5002
$ws.Range("M107").Value = -1364.5002

$ws.Range("H122").Value = 5483.136
$ws.Range("I122").Value = 3979.889
$ws.Range("K122").Value = 11939.667
$ws.Range("M122").Value = -9489.667000000001

$ws.Range("H125").Value = 134309.2
$ws.Range("J125").Value = 134309.2
$ws.Range("L125").Value = 134309.2
$ws.Range("N125").Value = -144149.2

$ws.Range("H132").Value = 4557.4736
$ws.Range("I132").Value = 4437.1562
$ws.Range("K132").Value = 13311.4686
$ws.Range("M132").Value = -10781.4686

$ws.Range("H137").Value = 61747
$ws.Range("J137").Value = 61747
$ws.Range("L137").Value = 61747
$ws.Range("N137").Value = -71947
